$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The evaluation window rolled forward by one quarter: for each row, the
# error statistics (columns B:G) that used to belong to the next quarter
# now belong to this quarter, and a brand-new set of statistics for the
# newest quarter is written into the last row (row 11). Column A (the
# Q0..Q9 labels) is untouched.

$newValues = @(
    @(0.01177902937386835, 0.3233677337305382, 0.1387525778210147, 0.3724950708680783, 0.3863630867095333, 14),
    @(-0.008234150994489252, 0.3656812685479343, 0.1736781030421116, 0.4167470492302394, 0.4336794057244024, 13),
    @(0.04114622104959433, 0.2438785208639479, 0.09301430872978532, 0.3049824728239072, 0.315631483433317, 12),
    @(-0.004251297071472607, 0.3196557001756482, 0.1664324393548997, 0.4079613209054257, 0.4278502102499434, 11),
    @(-0.03152289704198964, 0.3137577720403564, 0.1200567451426497, 0.3464920563918454, 0.3637200573468097, 10),
    @(-0.05669611638992965, 0.309241221828214, 0.12362493933599, 0.3516033835673229, 0.3680513563860018, 9),
    @(-0.04410238231315148, 0.2487243691024244, 0.1063280554679864, 0.3260798299005727, 0.3453909343541705, 8),
    @(-0.07289500678705241, 0.3460499651212735, 0.2127274776992719, 0.4612238910759848, 0.4919174337551497, 7),
    @(-0.05253816173664939, 0.2093953811075043, 0.05818608614024737, 0.2412179225104291, 0.2578972532655152, 6),
    @(-0.009352853557239749, 0.4545107205050175, 0.223274494992554, 0.4725193064759936, 0.5281891459539978, 5)
)

$row = 2
foreach ($rowValues in $newValues) {
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
    $ws.Cells.Item($row, 5).Value = $rowValues[3]
    $ws.Cells.Item($row, 6).Value = $rowValues[4]
    $ws.Cells.Item($row, 7).Value = $rowValues[5]
    $row++
}
